$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 119, shifting the existing rows (119..166) down to (120..167)
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record's data
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44559
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = 100112052
$ws.Cells.Item(119, 7).Value = "Albahaca"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 65
$ws.Cells.Item(119, 11).Value = 5000
$ws.Cells.Item(119, 12).Value = 5000
$ws.Cells.Item(119, 13).Value = 5000
$ws.Cells.Item(119, 14).Value = "$/paquete"
$ws.Cells.Item(119, 15).Value = "Región del Maule"
$ws.Cells.Item(119, 16).Value = 5000
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"
